# Applies the target edit: adds a new "Queries" worksheet after the
# existing report sheets, populated with report-query metadata (headers,
# agent/report descriptors, the before/after date bounds used by the
# report, and the three SQL queries driving it), then makes it the active
# sheet. Matches the diff against SkillHistoricalReportData.xlsx.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet right after the last existing sheet ---------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Queries"

# --- Header row (row 1) ----------------------------------------------------
$ws.Range("A1").Value = "Report Channel"
$ws.Range("B1").Value = "Report Name"
$ws.Range("C1").Value = "Report Type"
$ws.Range("D1").Value = "Start Date"
$ws.Range("E1").Value = "End Date"
$ws.Range("A1:E1").NumberFormat = "@"

# F1/G1/H1 (Query / QueryDrillGridOne / QueryDrillGridTwo) and the matching
# SQL text in row 2 must be written in this exact order so the new entries
# appended to the shared-string table line up with the source workbook.
$ws.Range("F1").Value = "Query"
$sql1 = @'
SELECT CONVERT(decimal(10,2), (100 * (sum(isnull(CallsHandledWithinSLAThreshold,0)))/
 (CAST(ISNULL((CASE WHEN SUM(isnull(PassedCalls,0))+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) = 0 THEN 1
 else SUM(PassedCalls)+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) end),1) AS float)))) AS [Service Level],
 SkillName as [Skill Name],sum([FlowIn]) AS [Flow In],sum([FlowOut]) AS [Flow Out],[SkillId] as [Skill ID],
 [dbo].[SECONDSTOhhmmss](sum(TotalStaffedTIme)/nullif(sum(TotalStaffedAgents),0)) AS [Avg Staff Time],
 [dbo].[SECONDSTOhhmmss](sum(TotalAbandTime)/nullif((sum(AbandCalls)+sum([AcdCalls])),0)) AS [Avg Aband Time],
 sum([AbandCalls]) AS [Aband Calls],
 [dbo].[SECONDSTOhhmmss](sum(SpeedOfAnswer)/nullif(sum([AcdCalls]),0)) AS [Avg Speed Answer],
 [dbo].[SECONDSTOhhmmss](sum([TotalAfterCallTime])) AS [Total After Call Time],
 [dbo].[SECONDSTOhhmmss](sum(TotalTalkTime)/nullif(sum([AcdCalls]),0)) AS [Avg Talk Time],
 sum([AcdCalls])  AS [Total Interaction],
 [dbo].[SECONDSTOhhmmss](sum(TotalAuxTime)) AS [Total Aux Time]
 from [OCM_SkillHistoricalReport]  WITH (NOLOCK)
 WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate' 
 GROUP BY [SkillId],[SkillName]
 ORDER BY [SkillName]
'@
$ws.Range("F2").Value = $sql1

$ws.Range("G1").Value = "QueryDrillGridOne"
$sql2 = @'
SELECT CONVERT(decimal(10,2), (100 * (sum(isnull(CallsHandledWithinSLAThreshold,0)))/
(CAST(ISNULL((CASE WHEN SUM(isnull(PassedCalls,0))+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) = 0 THEN 1
else SUM(PassedCalls)+sum(isnull(CallsAbandonedAfterSLAThreshold,0)) end),1) AS float))))  AS [Service Level],
sum([FlowIn]) AS [Flow In],sum([FlowOut]) AS [Flow Out],
[dbo].[SECONDSTOhhmmss](sum(TotalAbandTime)/nullif((sum(AbandCalls)+sum([AcdCalls])),0)) AS [Avg Aband Time],
sum([AbandCalls]) AS [Aband Calls],
Dateint AS [Date],[dbo].[SECONDSTOhhmmss](sum(SpeedOfAnswer)/nullif(sum([AcdCalls]),0)) AS [Avg Speed Answer], 
[dbo].[SECONDSTOhhmmss](sum([TotalAfterCallTime])) AS [Total After Call Time],
[dbo].[SECONDSTOhhmmss](sum(TotalTalkTime)/nullif(sum([AcdCalls]),0)) AS [Avg Talk Time],
sum([AcdCalls]) AS [Total Interaction],
[dbo].[SECONDSTOhhmmss](sum(TotalStaffedTIme)/nullif(sum(TotalStaffedAgents),0)) AS [Avg Staff],
[dbo].[SECONDSTOhhmmss](sum(TotalAuxTime)) AS [Total Aux Time]
from [OCM_SkillHistoricalReport] WITH (NOLOCK)
WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate' and  SkillId like 'SkillIdCapturedFromUI'
GROUP BY [Dateint],[SkillId],[SkillName] ORDER BY [Dateint] ASC
'@
$ws.Range("G2").Value = $sql2

$ws.Range("H1").Value = "QueryDrillGridTwo"
$sql3 = @'
SELECT [ServiceLevel] AS [Service Level],[FlowIn] AS [Flow In],[FlowOut] AS [Flow Out],
[dbo].[SECONDSTOhhmmss](TotalStaffedTIme/nullif(TotalStaffedAgents,0)) AS [Avg Staff Time],
[dbo].[SECONDSTOhhmmss](TotalAbandTime/nullif((AbandCalls+[AcdCalls]),0)) AS [Avg Aband Time],
[AbandCalls] AS [Aband Calls],
[dbo].[SECONDSTOhhmmss](SpeedOfAnswer/nullif([AcdCalls],0)) AS [Avg Speed Answer],
[dbo].[SECONDSTOhhmmss]([TotalAfterCallTime]) AS [Total After Call Time],
[dbo].[SECONDSTOhhmmss](TotalTalkTime/nullif([AcdCalls],0)) AS [Avg Talk Time],
[AcdCalls] AS [Total Interaction],[dbo].[SECONDSTOhhmmss](TotalAuxTime) AS [Total Aux Time],[Interval]
FROM [OCM_SkillHistoricalReport] WITH (NOLOCK) WHERE [ReportDateTime]>='ReportBeforeDate' AND [ReportDateTime]<='ReportAfterDate' AND 
[SkillId] like 'SkillIdCapturedFromUI' and [Interval] like '%' 
ORDER BY [intvl] ASC
'@
$ws.Range("H2").Value = $sql3

$ws.Range("G1").NumberFormat = "@"
$ws.Range("H1").NumberFormat = "@"

# --- Data row (row 2) -------------------------------------------------------
$ws.Range("A2").Value = "Agent"
$ws.Range("B2").Value = "OCM Skill Historical Report"
$ws.Range("C2").Value = "Date Range"

# Report date bounds, quote-prefixed so Excel keeps them as literal text
# instead of auto-converting them to date serials.
$ws.Range("D2").Value = "'08-04-2020 00:00:00"
$ws.Range("E2").Value = "'22-04-2020 00:00:00"

# Wrap the long SQL text cells and expand the row to fit the content.
$ws.Range("F2:H2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 409.5

# --- Column sizing -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Columns.Item(2).ColumnWidth = 24.14
$ws.Columns.Item(3).ColumnWidth = 10.14
$ws.Columns.Item(4).ColumnWidth = 17
$ws.Columns.Item(5).ColumnWidth = 17
$ws.Columns.Item(6).ColumnWidth = 20.43
$ws.Columns.Item(7).ColumnWidth = 15.57
$ws.Columns.Item(8).ColumnWidth = 15.29

# --- View state: make the new sheet active with E2 selected -------------
$ws.Range("E2").Select()
$ws.Activate()

# Use automatic recalculation (matches the workbook's post-edit calcPr).
$excel.Calculation = -4105

Write-Host "Queries sheet added"
